# Update "想去人数" (want-to-go count) figures in column F for the
# "展览" (Exhibitions) and "全部类型" (All types) sheets.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 386
$wsExhibit.Range("F3").Value = 0
$wsExhibit.Range("F5").Value = 0
$wsExhibit.Range("F6").Value = 23
$wsExhibit.Range("F8").Value = 0
$wsExhibit.Range("F10").Value = 481

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 386
$wsAll.Range("F3").Value = 0
$wsAll.Range("F4").Value = 1619
$wsAll.Range("F5").Value = 0
$wsAll.Range("F6").Value = 0
$wsAll.Range("F10").Value = 481
